# horarios-141 arribos update (2026-01-09 run)
# Refreshes the ETA/BANDERA/MIN/ESTADO rows on TODOS, 215 and COMBINADAS with the
# latest arrival snapshot, and trims each sheet to its new (shorter) row count.

$wb = $excel.ActiveWorkbook

# New data for TODOS / COMBINADAS (29 data rows: A2:D29)
$todosData = @(
    @('17:42', '27_EL RETIRO', 0, '🚌'),
    @('17:45', '10_OLMOS', 3, '🚌'),
    @('17:51', '16_P MOR-167 Y 521', 9, '🚌'),
    @('17:52', '81_EL PELIGRO', 10, '📅'),
    @('17:53', '11_ETCHEVERRY', 11, '📅'),
    @('18:01', '16_SANTA ANA', 19, '🚌'),
    @('18:04', '17_ROMERO', 22, '🚌'),
    @('18:04', '215C_LA PLATA', 22, '🚌'),
    @('18:05', '23_HERNANDEZ', 23, '🚌'),
    @('18:09', '14_ABASTO', 27, '🚌'),
    @('18:11', '16_SANTA ANA', 29, '🚌'),
    @('18:16', '10_OLMOS', 34, '🚌'),
    @('18:16', '15_ABASTO', 34, '🚌'),
    @('18:21', '26_HERNANDEZ', 39, '🚌'),
    @('18:25', '14_ABASTO', 43, '🚌'),
    @('18:28', '215C_EL PATO', 46, '🚌'),
    @('18:32', '11X44_ETCHEVERRY', 50, '🚌'),
    @('18:35', '23_HERNANDEZ', 53, '🚌'),
    @('18:40', '15_ABASTO', 58, '🚌'),
    @('18:48', '14X44_ABASTO', 66, '🚌'),
    @('18:52', '215A_LA PLATA', 70, '🚌'),
    @('18:59', '215A_EL PATO', 77, '📅'),
    @('19:04', '215B_LP-P MOR-1 Y 57', 82, '📅'),
    @('19:05', '11_ETCHEVERRY', 83, '📅'),
    @('19:11', '16_P MOR-SANTA ANA', 89, '🚌'),
    @('19:17', '27_EL RETIRO', 95, '🚌'),
    @('19:30', '225_GOMEZ', 108, '📅'),
    @('19:40', '215C_EL PATO', 118, '🚌')
)

# New data for 215 (7 data rows: A2:D7) - the "215*" BANDERA subset of TODOS
$data215 = @(
    @('18:04', '215C_LA PLATA', 22, '🚌'),
    @('18:28', '215C_EL PATO', 46, '🚌'),
    @('18:52', '215A_LA PLATA', 70, '🚌'),
    @('18:59', '215A_EL PATO', 77, '📅'),
    @('19:04', '215B_LP-P MOR-1 Y 57', 82, '📅'),
    @('19:40', '215C_EL PATO', 118, '🚌')
)

function Update-Sheet($SheetName, $Rows) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Overwrite existing header-adjacent rows with the refreshed values.
    $r = 2
    foreach ($row in $Rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $r++
    }

    # Remove any now-stale trailing rows left over from the longer, previous table.
    $lastNewRow = $Rows.Count + 1
    $oldLastRow = $ws.UsedRange.Rows.Count
    for ($row = $oldLastRow; $row -gt $lastNewRow; $row--) {
        $ws.Rows.Item($row).Delete()
    }
}

Update-Sheet "TODOS" $todosData
Update-Sheet "215" $data215
Update-Sheet "COMBINADAS" $todosData
